$wb = $excel.ActiveWorkbook

# Sheet "展览" - update 想去人数 (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 56
$ws1.Range("F3").Value = 339
$ws1.Range("F4").Value = 23
$ws1.Range("F5").Value = 91

# Sheet "全部类型" - same updates mirrored from "展览"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 56
$ws4.Range("F3").Value = 339
$ws4.Range("F4").Value = 23
$ws4.Range("F5").Value = 91
